# Update "gh-pages" generated numbers (想去人数 / 最低票价) per latest scrape.
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 767
$ws1.Range("F4").Value  = 1538
$ws1.Range("F5").Value  = 236
$ws1.Range("F7").Value  = 165
$ws1.Range("F8").Value  = 6322
$ws1.Range("F12").Value = 5382
$ws1.Range("F15").Value = 1203
$ws1.Range("F16").Value = 3
$ws1.Range("F18").Value = 367
$ws1.Range("F21").Value = 308
$ws1.Range("F22").Value = 30
$ws1.Range("F24").Value = 3849
$ws1.Range("F25").Value = 164

# ---- Sheet: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 96
$ws2.Range("G2").Value = 299

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 96
$ws4.Range("G2").Value  = 299
$ws4.Range("F4").Value  = 767
$ws4.Range("F5").Value  = 1538
$ws4.Range("F6").Value  = 236
$ws4.Range("F8").Value  = 165
$ws4.Range("F9").Value  = 6322
$ws4.Range("F13").Value = 5382
$ws4.Range("F16").Value = 1203
$ws4.Range("F17").Value = 3
$ws4.Range("F19").Value = 367
$ws4.Range("F22").Value = 308
$ws4.Range("F23").Value = 30
$ws4.Range("F25").Value = 3850
$ws4.Range("F27").Value = 164
